# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.712.76"
$ws.Range("E2").Value = "  -2.31%  "

$ws.Range("D3").Value = "2.342.36"
$ws.Range("E3").Value = "  -2.54%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").Value = "'499.07"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "'128.01"
$ws.Range("E6").Value = "  -3.82%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("D9").Value = "2.349.55"
$ws.Range("E9").Value = "  -2.61%  "

$ws.Range("D10").Value = "'0.0975"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'4.84"
$ws.Range("E12").Value = "  +5.45%  "

$ws.Range("D13").Value = "'0.320"
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("D14").Value = "2.776.32"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").Value = "55.736.30"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").Value = "'21.47"
$ws.Range("E16").Value = "  -1.74%  "

$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("D18").Value = "2.339.50"
$ws.Range("E18").Value = "  -3.96%  "

$ws.Range("D19").Value = "'9.89"
$ws.Range("E19").Value = "  -3.74%  "

$ws.Range("D20").Value = "'307.39"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").Value = "'3.98"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  -1.37%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "'65.13"
$ws.Range("E24").Value = "  -2.88%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("D27").Value = "'0.146"
$ws.Range("E27").Value = "  -4.31%  "

$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = "  -4.84%  "

$ws.Range("D29").Value = "'172.71"
$ws.Range("E29").Value = "  -1.68%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  -2.95%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0699"
$ws.Range("E31").Value = "  -3.95%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.77"
$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").Value = "'1.06"
$ws.Range("E35").Value = "  -6.09%  "

$ws.Range("D36").Value = "'17.43"
$ws.Range("E36").Value = "  -3.06%  "

$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = "  -2.57%  "

$ws.Range("D38").Value = "'3.62"
$ws.Range("E38").Value = "  -5.70%  "

$ws.Range("D39").Value = "'0.817"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").Value = "'36.20"
$ws.Range("E40").Value = "  -1.81%  "

$ws.Range("D41").Value = "'1.37"
$ws.Range("E41").Value = "  -5.31%  "

$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("D43").Value = "'125.97"
$ws.Range("E43").Value = "  -5.50%  "

$ws.Range("D44").Value = "'4.67"
$ws.Range("E44").Value = "  -4.52%  "

$ws.Range("D45").Value = "'0.557"
$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("D46").Value = "'0.0891"
$ws.Range("E46").Value = "  -2.44%  "

$ws.Range("D47").Value = "'234.42"
$ws.Range("E47").Value = "  -6.92%  "

$ws.Range("D48").Value = "'0.0476"
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("D49").Value = "'0.0204"
$ws.Range("E49").Value = "  -3.72%  "

$ws.Range("D50").Value = "'16.55"
$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("E51").Value = "  +0.28%  "

